$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 116.9
$ws.Range("J9").Value = 164.66667
$ws.Range("L9").Value = 164.66667
$ws.Range("N9").Value = -502.66667
$ws.Range("H12").Value = 111.111115
$ws.Range("J12").Value = 150
$ws.Range("L12").Value = 150
$ws.Range("N12").Value = -490
$ws.Range("H38").Value = 1854.6842
$ws.Range("I38").Value = 121.42857
$ws.Range("J38").Value = 2865.75
$ws.Range("K38").Value = 364.28571
$ws.Range("L38").Value = 8597.25
$ws.Range("M38").Value = 7.714290000000005
$ws.Range("N38").Value = -9341.25
$ws.Range("H40").Value = 798
$ws.Range("I40").Value = 798
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 798
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -623
$ws.Range("N40").ClearContents()
$ws.Range("H62").Value = 159099520
$ws.Range("I62").Value = 71441830
$ws.Range("J62").Value = 312500500
$ws.Range("K62").Value = 71441830
$ws.Range("L62").Value = 312500500
$ws.Range("M62").Value = -71441206
$ws.Range("N62").Value = -312501748
$ws.Range("H65").Value = 159099520
$ws.Range("I65").Value = 71441830
$ws.Range("J65").Value = 312500500
$ws.Range("K65").Value = 357209150
$ws.Range("L65").Value = 1562502500
$ws.Range("M65").Value = -357206030
$ws.Range("N65").Value = -1562508740
$ws.Range("H98").Value = 20929934
$ws.Range("I98").Value = 7408795.5
$ws.Range("K98").Value = 7408795.5
$ws.Range("M98").Value = -7407297.5
$ws.Range("H107").Value = 2159.5557
$ws.Range("I107").Value = 3430
$ws.Range("J107").Value = 571.5
$ws.Range("K107").Value = 3430
$ws.Range("L107").Value = 571.5
$ws.Range("M107").Value = -1510
$ws.Range("N107").Value = -4411.5
$ws.Range("H122").Value = 20929934
$ws.Range("I122").Value = 7408795.5
$ws.Range("K122").Value = 22226386.5
$ws.Range("M122").Value = -22223936.5
$ws.Range("H132").Value = 2471188.8
$ws.Range("I132").Value = 1939.079
$ws.Range("J132").Value = 15875687
$ws.Range("K132").Value = 5817.237
$ws.Range("L132").Value = 47627061
$ws.Range("M132").Value = -3287.237
$ws.Range("N132").Value = -47632121

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 8587020
$ws.Range("I61").Value = 7921734.5
$ws.Range("J61").Value = 9806709
$ws.Range("K61").Value = 7921734.5
$ws.Range("L61").Value = 9806709
$ws.Range("M61").Value = -7921522.5
$ws.Range("N61").Value = -9807133
$ws.Range("H122").Value = 2825.4285
$ws.Range("I122").Value = 2916.7144
$ws.Range("J122").Value = 2642.8572
$ws.Range("K122").Value = 8750.143199999999
$ws.Range("L122").Value = 7928.571599999999
$ws.Range("M122").Value = -6300.143199999999
$ws.Range("N122").Value = -12828.5716
$ws.Range("H132").Value = 20560236
$ws.Range("I132").Value = 18061384
$ws.Range("J132").Value = 25557936
$ws.Range("K132").Value = 54184152
$ws.Range("L132").Value = 76673808
$ws.Range("M132").Value = -54181622
$ws.Range("N132").Value = -76678868
$ws.Range("H136").Value = 8587020
$ws.Range("I136").Value = 7921734.5
$ws.Range("J136").Value = 9806709
$ws.Range("K136").Value = 23765203.5
$ws.Range("L136").Value = 29420127
$ws.Range("M136").Value = -23762653.5
$ws.Range("N136").Value = -29425227

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H134").Value = 24065068
$ws.Range("I134").Value = 29412560
$ws.Range("J134").Value = 5883595.5
$ws.Range("K134").Value = 88237680
$ws.Range("L134").Value = 17650786.5
$ws.Range("M134").Value = -88235145
$ws.Range("N134").Value = -17655856.5

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 289.8095
$ws.Range("I22").Value = 234.3077
$ws.Range("J22").Value = 380
$ws.Range("K22").Value = 234.3077
$ws.Range("L22").Value = 380
$ws.Range("M22").Value = 115.6923
$ws.Range("N22").Value = -1080
$ws.Range("H50").Value = 19992
$ws.Range("J50").Value = 19992
$ws.Range("L50").Value = 19992
$ws.Range("N50").Value = -21242
$ws.Range("H51").Value = 19565.666
$ws.Range("J51").Value = 19565.666
$ws.Range("L51").Value = 19565.666
$ws.Range("N51").Value = -21037.666
$ws.Range("H58").Value = 1178221.5
$ws.Range("I58").Value = 1731028.4
$ws.Range("J58").Value = 3507
$ws.Range("K58").Value = 1731028.4
$ws.Range("L58").Value = 3507
$ws.Range("M58").Value = -1730825.4
$ws.Range("N58").Value = -3913
$ws.Range("H60").Value = 7493.5
$ws.Range("I60").Value = 4993.2
$ws.Range("J60").Value = 19995
$ws.Range("K60").Value = 4993.2
$ws.Range("L60").Value = 19995
$ws.Range("M60").Value = -4482.2
$ws.Range("N60").Value = -21017
$ws.Range("H61").Value = 19565.666
$ws.Range("J61").Value = 19565.666
$ws.Range("L61").Value = 19565.666
$ws.Range("N61").Value = -20261.666
$ws.Range("H122").Value = 15585.333
$ws.Range("I122").Value = 18202.4
$ws.Range("K122").Value = 54607.2
$ws.Range("M122").Value = -52157.2
$ws.Range("H134").Value = 1178559.5
$ws.Range("I134").Value = 1563.2593
$ws.Range("J134").Value = 5718402
$ws.Range("K134").Value = 4689.7779
$ws.Range("L134").Value = 17155206
$ws.Range("M134").Value = -2154.7779
$ws.Range("N134").Value = -17160276
$ws.Range("H136").Value = 1178221.5
$ws.Range("I136").Value = 1731028.4
$ws.Range("J136").Value = 3507
$ws.Range("K136").Value = 5193085.199999999
$ws.Range("L136").Value = 10521
$ws.Range("M136").Value = -5190535.199999999
$ws.Range("N136").Value = -15621

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 103.51613
$ws.Range("I2").Value = 125.07692
$ws.Range("J2").Value = 87.94444
$ws.Range("K2").Value = 750.4615200000001
$ws.Range("L2").Value = 527.66664
$ws.Range("M2").Value = -637.4615200000001
$ws.Range("N2").Value = -753.66664
$ws.Range("H63").Value = 3048.3333
$ws.Range("I63").Value = 697.5
$ws.Range("J63").Value = 7750
$ws.Range("K63").Value = 2092.5
$ws.Range("L63").Value = 23250
$ws.Range("M63").Value = -1343.5
$ws.Range("N63").Value = -24748
$ws.Range("H64").Value = 4459.3335
$ws.Range("I64").Value = 12
$ws.Range("J64").Value = 4863.636
$ws.Range("K64").Value = 36
$ws.Range("L64").Value = 14590.908
$ws.Range("M64").Value = 234
$ws.Range("N64").Value = -15130.908
$ws.Range("H66").Value = 3048.3333
$ws.Range("I66").Value = 697.5
$ws.Range("J66").Value = 7750
$ws.Range("K66").Value = 6277.5
$ws.Range("L66").Value = 69750
$ws.Range("M66").Value = -2533.5
$ws.Range("N66").Value = -77238
$ws.Range("H67").Value = 4459.3335
$ws.Range("I67").Value = 12
$ws.Range("J67").Value = 4863.636
$ws.Range("K67").Value = 36
$ws.Range("L67").Value = 14590.908
$ws.Range("M67").Value = 900
$ws.Range("N67").Value = -16462.908
$ws.Range("H68").Value = 2763.1296
$ws.Range("I68").Value = 548.10345
$ws.Range("J68").Value = 5332.56
$ws.Range("K68").Value = 1644.31035
$ws.Range("L68").Value = 15997.68
$ws.Range("M68").Value = -833.3103499999997
$ws.Range("N68").Value = -17619.68
$ws.Range("H71").Value = 2763.1296
$ws.Range("I71").Value = 548.10345
$ws.Range("J71").Value = 5332.56
$ws.Range("K71").Value = 4932.931049999999
$ws.Range("L71").Value = 47993.04
$ws.Range("M71").Value = -876.9310499999992
$ws.Range("N71").Value = -56105.04
$ws.Range("H75").Value = 1166.6666
$ws.Range("J75").Value = 1500
$ws.Range("L75").Value = 4500
$ws.Range("N75").Value = -6496
$ws.Range("H78").Value = 1166.6666
$ws.Range("J78").Value = 1500
$ws.Range("L78").Value = 13500
$ws.Range("N78").Value = -23484
$ws.Range("H112").Value = 3974.1538
$ws.Range("J112").Value = 5095.5557
$ws.Range("L112").Value = 15286.6671
$ws.Range("N112").Value = -17502.6671
$ws.Range("H119").Value = 8192.9
$ws.Range("I119").Value = 964.5
$ws.Range("J119").Value = 10000
$ws.Range("K119").Value = 2893.5
$ws.Range("L119").Value = 30000
$ws.Range("M119").Value = 1944.5
$ws.Range("N119").Value = -39676
$ws.Range("H129").Value = 34633996
$ws.Range("I129").Value = 90911750
$ws.Range("J129").Value = 6495117.5
$ws.Range("K129").Value = 272735250
$ws.Range("L129").Value = 19485352.5
$ws.Range("M129").Value = -272730250
$ws.Range("N129").Value = -19495352.5
$ws.Range("H131").Value = 7829048.5
$ws.Range("J131").Value = 17634.35
$ws.Range("L131").Value = 52903.05
$ws.Range("N131").Value = -62983.05

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 9806604
$ws.Range("I122").Value = 2383.3333
$ws.Range("K122").Value = 7149.999899999999
$ws.Range("M122").Value = -4699.999899999999
$ws.Range("H132").Value = 19324856
$ws.Range("I132").Value = 15333976
$ws.Range("J132").Value = 25976324
$ws.Range("K132").Value = 46001928
$ws.Range("L132").Value = 77928972
$ws.Range("M132").Value = -45999398
$ws.Range("N132").Value = -77934032

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 166671660
$ws.Range("I46").Value = 10000
$ws.Range("J46").Value = 250002500
$ws.Range("K46").Value = 10000
$ws.Range("L46").Value = 250002500
$ws.Range("M46").Value = -9812
$ws.Range("N46").Value = -250002876
$ws.Range("H93").Value = 26092.857
$ws.Range("I93").Value = 14679.556
$ws.Range("J93").Value = 46636.8
$ws.Range("K93").Value = 14679.556
$ws.Range("L93").Value = 46636.8
$ws.Range("M93").Value = -13431.556
$ws.Range("N93").Value = -49132.8
